$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.391.97"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "  +1.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.939.54"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.52"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  +3.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.99"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "  +2.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.734"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("E10").Value = "  +3.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000356"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = "  +5.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.16"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.76"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = "  +3.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.566.56"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.75"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = "  -2.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.958.40"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.02"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.14"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.507.32"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.27"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.48"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "  +4.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.83"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.80"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.44"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = "  +13.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.07"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = "  +15.61%  "

$ws.Range("E27").Value = "  +2.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.92"
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = "  +0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.88"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = "  +2.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "715.10"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.69"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = "  +3.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0924"
$ws.Range("D34").ClearFormats()

$ws.Range("E34").Value = "  +15.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.08"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.16"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = "  +14.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "61.15"
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = "  +5.78%  "

$ws.Range("E38").Value = "  -3.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.401"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = "  +19.26%  "

$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0483"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = "  +1.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  +14.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.17"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "  +4.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.96"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = "  +5.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.143"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("B47").Value = "BabyDogeCoin"

$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0360"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "  +40.53%  "

$ws.Range("B48").Value = "LidoDAOToken"

$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.43"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.89"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("B51").Value = "FLOKI"

$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000270"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").Value = "  +46.02%  "
